$d = $word.ActiveDocument

# 1. Update the date/time stamp paragraph.
$d.Content.Find.Execute(
    "May  31, 2021 (10:16:19 PM)", $true, $false, $false, $false, $false,
    $true, 1, $false, "June   1, 2021 (01:00:30 AM)", 2)

# 2. Merge the three-sentence intro paragraph into one run, dropping the
#    middle sentence ("We will use the example shown in lecture.").
$d.Content.Find.Execute(
    "This lab will guide you in your first manipulation of a programmer-defined class. We will use the example shown in lecture. The last part is challenging; therefore, we provide a possible solution at the end of the page, but make sure you try to solve it by yourself beforehand.",
    $true, $false, $false, $false, $false,
    $true, 1, $false,
    "This lab will guide you in your first manipulation of a programmer-defined class. The last part is challenging; therefore, we provide a possible solution at the end of the page, but make sure you try to solve it by yourself beforehand.",
    2)
